$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.121.73'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '2.420.68'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '569.62'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('D6').Value = '139.45'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D9').Value = '2.406.33'
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').Value = '2.857.52'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = '61.030.30'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = '2.415.15'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '7.86'
$ws.Range('E19').Value = '  +9.08%  '
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('D21').Value = '323.32'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('E23').Value = '  +1.89%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '1.83'
$ws.Range('E25').Value = '  -4.25%  '
$ws.Range('D26').Value = '64.73'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '582.84'
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('E28').Value = '  -9.77%  '
$ws.Range('D29').Value = '2.536.22'
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('D30').Value = '0.0₃0933'
$ws.Range('E30').Value = '  -3.10%  '
$ws.Range('E32').Value = '  -5.17%  '
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '4.62'
$ws.Range('E37').Value = '  -5.19%  '
$ws.Range('D38').Value = '151.85'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('D40').Value = '18.23'
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('E41').Value = '  -2.40%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('D44').Value = '41.09'
$ws.Range('E44').Value = '  -4.76%  '
$ws.Range('D45').Value = '2.34'
$ws.Range('E45').Value = '  -5.98%  '
$ws.Range('D46').Value = '142.91'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('E48').Value = '  -2.94%  '
$ws.Range('D49').Value = '0.587'
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').Value = '19.51'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('E51').Value = '  -3.22%  '
